$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "Valentina Perghem "
$ws.Range("B16").Value = "Alessandro Festi | La Contea FC"
$ws.Range("C16").Value = "Edoardo Pomarolli | Modium"
$ws.Range("D16").Value = "Luca Frasca | Clitoriders"
$ws.Range("E16").Value = "Niccolò Orsi | SBARX"
$ws.Range("F16").Value = "Andrea Menolli | SdrumALA"
